$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: header date (Sunday 03-03-2024) formatted as a date
$ws.Range("E1").Value = 45354
$ws.Range("E1").NumberFormat = "mm-dd-yy"

# New column E data values
$ws.Range("E2").Value = 1027
$ws.Range("E3").Value = 998
$ws.Range("E4").Value = 1500
$ws.Range("E5").Value = 2919
$ws.Range("E6").Value = 95

# Column E width (best-fit sized like the source diff)
$ws.Columns.Item(5).ColumnWidth = 9.43

# Match the saved selection state from the edit (active cell E3)
$ws.Range("E3").Select()
